$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (45179 -> 45180) for every data row (rows 2 through 471).
$ws.Range("C2:C471").Value = 45180
